$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# These numeric-looking values are stored as text in the workbook, so
# prefix with an apostrophe to keep them text (matching the existing
# cell type) instead of letting Excel coerce them to numbers.

# Enterprises density (per 1000 people) row 11: Micro/SMEs/MSMEs
$ws.Range("B11").Value = "'0.08"
$ws.Range("C11").Value = "'3.68"
$ws.Range("D11").Value = "'3.76"

# Employment (% of total) row 12: Micro/SMEs/MSMEs
$ws.Range("B12").Value = "'7.78"
$ws.Range("C12").Value = "'38.97"
$ws.Range("D12").Value = "'46.75"

# Enterprises (% of total) row 14: Micro/SMEs/MSMEs
$ws.Range("B14").Value = "'2.01"
$ws.Range("C14").Value = "'95.72"
$ws.Range("D14").Value = "'97.73"
